$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: add columns P and Q, copying O1 style (bold/border/center), with values 14 and 15 ---
$ws.Range("O1").Copy() | Out-Null
$ws.Range("P1:Q1").PasteSpecial(-4122) | Out-Null
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# --- Data rows 2-25: recomputed stats -> update B, C, E, (F minor), G, H, I;
#     O zeroed out, new P column is 0, new Q column carries the old O value ---
# row 2
$ws.Range("B2").Value = 3.45706446842388
$ws.Range("C2").Value = 1.007463387640854
$ws.Range("E2").Value = 1.343001604769114
$ws.Range("G2").Value = 0.0007919315264227313
$ws.Range("H2").Value = 0.009976046992029766
$ws.Range("I2").Value = 0.002654056212566314
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 0.9101585796577893

# row 3
$ws.Range("B3").Value = 3.012171406964342
$ws.Range("C3").Value = 0.8887149943078327
$ws.Range("E3").Value = 1.169389773686817
$ws.Range("F3").Value = 1.26082040731464
$ws.Range("G3").Value = 0.0007954614514921511
$ws.Range("H3").Value = 0.006951618590052611
$ws.Range("I3").Value = 0.001410059181311318
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = 0
$ws.Range("Q3").Value = 0.8891057585376387

# row 4
$ws.Range("B4").Value = 2.738754624972501
$ws.Range("C4").Value = 0.8163771175811689
$ws.Range("E4").Value = 1.063116433664888
$ws.Range("G4").Value = 0.0007976930960263985
$ws.Range("H4").Value = 0.005315784855754235
$ws.Range("I4").Value = 0.0008979189225590822
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 0
$ws.Range("Q4").Value = 0.8782150267306292

# row 5
$ws.Range("B5").Value = 2.627221208681419
$ws.Range("C5").Value = 0.7886090506703738
$ws.Range("E5").Value = 1.019853648607011
$ws.Range("F5").Value = 1.090973950927975
$ws.Range("G5").Value = 0.0007986238881468814
$ws.Range("H5").Value = 0.004698928024502003
$ws.Range("I5").Value = 0.0008088876672678325
$ws.Range("O5").Value = 0
$ws.Range("P5").Value = 0
$ws.Range("Q5").Value = 0.8717184442198374

# row 6
$ws.Range("B6").Value = 2.608661551526438
$ws.Range("C6").Value = 0.7859419212789192
$ws.Range("E6").Value = 1.012653962182483
$ws.Range("F6").Value = 1.082851457025441
$ws.Range("G6").Value = 0.0007987853965448778
$ws.Range("H6").Value = 0.004597866120651273
$ws.Range("I6").Value = 0.0008797876039716712
$ws.Range("O6").Value = 0
$ws.Range("P6").Value = 0
$ws.Range("Q6").Value = 0.8675813516168773

# row 7
$ws.Range("B7").Value = 2.737157128669196
$ws.Range("C7").Value = 0.8213050690186208
$ws.Range("E7").Value = 1.062480831808543
$ws.Range("G7").Value = 0.0007977214043917563
$ws.Range("H7").Value = 0.005302660952479443
$ws.Range("I7").Value = 0.001105510563603751
$ws.Range("O7").Value = 0
$ws.Range("P7").Value = 0
$ws.Range("Q7").Value = 0.8696960583638429

# row 8
$ws.Range("B8").Value = 3.303572487780059
$ws.Range("C8").Value = 0.9734888820579783
$ws.Range("E8").Value = 1.28298990155524
$ws.Range("F8").Value = 1.390742831773508
$ws.Range("G8").Value = 0.0007931546222333261
$ws.Range("H8").Value = 0.008877294755860243
$ws.Range("I8").Value = 0.002415490825453226
$ws.Range("O8").Value = 0
$ws.Range("P8").Value = 0
$ws.Range("Q8").Value = 0.8913306168166599

# row 9
$ws.Range("B9").Value = 4.413697388953381
$ws.Range("C9").Value = 1.267469075783538
$ws.Range("E9").Value = 1.719131879611268
$ws.Range("F9").Value = 1.895061171842798
$ws.Range("G9").Value = 0.0007846681333266091
$ws.Range("H9").Value = 0.01786782060856185
$ws.Range("I9").Value = 0.006909485273287608
$ws.Range("O9").Value = 0
$ws.Range("P9").Value = 0
$ws.Range("Q9").Value = 0.9669023301322994

# row 10
$ws.Range("B10").Value = 5.229541537614182
$ws.Range("C10").Value = 1.489227463985628
$ws.Range("E10").Value = 1.941193055041964
$ws.Range("F10").Value = 2.275587181990815
$ws.Range("G10").Value = 0.0007789065123076898
$ws.Range("H10").Value = 0.02527337859572842
$ws.Range("I10").Value = 0.01186560577040563
$ws.Range("O10").Value = 0
$ws.Range("P10").Value = 0
$ws.Range("Q10").Value = 0.9923685502502906

# row 11
$ws.Range("B11").Value = 5.592035089196713
$ws.Range("C11").Value = 1.59189657791859
$ws.Range("E11").Value = 1.254694484650472
$ws.Range("F11").Value = 2.451628164639231
$ws.Range("G11").Value = 0.0007778393998584065
$ws.Range("H11").Value = 0.04071378888487232
$ws.Range("I11").Value = 0.01359257270796466
$ws.Range("O11").Value = 0
$ws.Range("P11").Value = 0
$ws.Range("Q11").Value = 0.6840291361303485

# row 12
$ws.Range("B12").Value = 5.726383921551019
$ws.Range("C12").Value = 1.623456135342963
$ws.Range("E12").Value = 0.7606207801202487
$ws.Range("G12").Value = 0.0007779067756840159
$ws.Range("H12").Value = 0.0766274760464043
$ws.Range("I12").Value = 0.01369507220244959
$ws.Range("O12").Value = 0
$ws.Range("P12").Value = 0
$ws.Range("Q12").Value = 0.4736524581405348

# row 13
$ws.Range("B13").Value = 5.69149969114909
$ws.Range("C13").Value = 1.610267069501333
$ws.Range("E13").Value = 0.3862190595630395
$ws.Range("F13").Value = 2.504295469238414
$ws.Range("G13").Value = 0.0007788722137611204
$ws.Range("H13").Value = 0.1296338476240209
$ws.Range("I13").Value = 0.01282651398845491
$ws.Range("O13").Value = 0
$ws.Range("P13").Value = 0
$ws.Range("Q13").Value = 0.3121078267026292

# row 14
$ws.Range("B14").Value = 5.590554727746905
$ws.Range("C14").Value = 1.582602072687394
$ws.Range("E14").Value = 0.1968044168632517
$ws.Range("G14").Value = 0.0007799330981433236
$ws.Range("H14").Value = 0.1771830353365829
$ws.Range("I14").Value = 0.01190201796664603
$ws.Range("O14").Value = 0
$ws.Range("P14").Value = 0
$ws.Range("Q14").Value = 0.2259690028565231

# row 15
$ws.Range("B15").Value = 5.529752952375816
$ws.Range("C15").Value = 1.568534547397974
$ws.Range("E15").Value = 0.1596373502882997
$ws.Range("G15").Value = 0.0007804154799479467
$ws.Range("H15").Value = 0.1890866575273265
$ws.Range("I15").Value = 0.01156379217992409
$ws.Range("O15").Value = 0
$ws.Range("P15").Value = 0
$ws.Range("Q15").Value = 0.208809031408407

# row 16
$ws.Range("B16").Value = 5.184712317784374
$ws.Range("C16").Value = 1.478695092013709
$ws.Range("E16").Value = 0.1555078664300424
$ws.Range("F16").Value = 2.264147245220144
$ws.Range("G16").Value = 0.0007826700725010032
$ws.Range("H16").Value = 0.1740758739912849
$ws.Range("I16").Value = 0.009660260516514363
$ws.Range("O16").Value = 0
$ws.Range("P16").Value = 0
$ws.Range("Q16").Value = 0.2274106021470175

# row 17
$ws.Range("B17").Value = 4.974309208171974
$ws.Range("C17").Value = 1.424871884375534
$ws.Range("E17").Value = 0.2391311429574401
$ws.Range("G17").Value = 0.000783856930264254
$ws.Range("H17").Value = 0.1354299262250294
$ws.Range("I17").Value = 0.0087354466932279
$ws.Range("O17").Value = 0
$ws.Range("P17").Value = 0
$ws.Range("Q17").Value = 0.2835228288309182

# row 18
$ws.Range("B18").Value = 4.855468628462802
$ws.Range("C18").Value = 1.39117380567194
$ws.Range("E18").Value = 0.4709913206318603
$ws.Range("F18").Value = 2.107028586650728
$ws.Range("G18").Value = 0.0007841387904206719
$ws.Range("H18").Value = 0.08325542121928464
$ws.Range("I18").Value = 0.008256881803588634
$ws.Range("O18").Value = 0
$ws.Range("P18").Value = 0
$ws.Range("Q18").Value = 0.4019045300307198

# row 19
$ws.Range("B19").Value = 4.818908855008374
$ws.Range("C19").Value = 1.387277969842842
$ws.Range("E19").Value = 0.9032220602888259
$ws.Range("F19").Value = 2.087706772602971
$ws.Range("G19").Value = 0.0007835304193727203
$ws.Range("H19").Value = 0.0406011370972692
$ws.Range("I19").Value = 0.008711353139416289
$ws.Range("O19").Value = 0
$ws.Range("P19").Value = 0
$ws.Range("Q19").Value = 0.5855600131483385

# row 20
$ws.Range("B20").Value = 5.014790210718388
$ws.Range("C20").Value = 1.447889257766235
$ws.Range("E20").Value = 1.878280411639778
$ws.Range("G20").Value = 0.0007804456614021833
$ws.Range("H20").Value = 0.02315411303990045
$ws.Range("I20").Value = 0.01107259061290655
$ws.Range("O20").Value = 0
$ws.Range("P20").Value = 0
$ws.Range("Q20").Value = 0.9575343383031196

# row 21
$ws.Range("B21").Value = 5.643542784591943
$ws.Range("C21").Value = 1.617633673222088
$ws.Range("E21").Value = 2.208159889130968
$ws.Range("G21").Value = 0.0007758007075602692
$ws.Range("H21").Value = 0.03052528701895785
$ws.Range("I21").Value = 0.01560759681299384
$ws.Range("O21").Value = 0
$ws.Range("P21").Value = 0
$ws.Range("Q21").Value = 1.049367805858594

# row 22
$ws.Range("B22").Value = 6.054981629084921
$ws.Range("C22").Value = 1.721897093229813
$ws.Range("E22").Value = 2.373131035669061
$ws.Range("F22").Value = 2.667416158031983
$ws.Range("G22").Value = 0.0007728912605559746
$ws.Range("H22").Value = 0.03536816661742437
$ws.Range("I22").Value = 0.01867232433822874
$ws.Range("O22").Value = 0
$ws.Range("P22").Value = 0
$ws.Range("Q22").Value = 1.103514759358632

# row 23
$ws.Range("B23").Value = 5.835425757555583
$ws.Range("C23").Value = 1.659583109929201
$ws.Range("E23").Value = 2.28502552055393
$ws.Range("G23").Value = 0.0007744277975906085
$ws.Range("H23").Value = 0.03275557330050738
$ws.Range("I23").Value = 0.01675250726775968
$ws.Range("O23").Value = 0
$ws.Range("P23").Value = 0
$ws.Range("Q23").Value = 1.083991483364542

# row 24
$ws.Range("B24").Value = 5.005576534307068
$ws.Range("C24").Value = 1.435901699072417
$ws.Range("E24").Value = 1.953721628580638
$ws.Range("G24").Value = 0.0007803510040827074
$ws.Range("H24").Value = 0.02362040906012486
$ws.Range("I24").Value = 0.01071079612140746
$ws.Range("O24").Value = 0
$ws.Range("P24").Value = 0
$ws.Range("Q24").Value = 1.002213127649725

# row 25
$ws.Range("B25").Value = 4.113089848590278
$ws.Range("C25").Value = 1.197150090379012
$ws.Range("E25").Value = 1.600502819685076
$ws.Range("F25").Value = 1.757113814925091
$ws.Range("G25").Value = 0.0007869364568551807
$ws.Range("H25").Value = 0.01518161825300501
$ws.Range("I25").Value = 0.005770901077604584
$ws.Range("O25").Value = 0
$ws.Range("P25").Value = 0
$ws.Range("Q25").Value = 0.9291665970757776

